# "Added getTarget to user rather than game."
#
# The GET sheet listed `/game/target` (and its sibling `/game/hugs`) in the
# "game" section. This moves that endpoint to the "user" section as
# `/user/target`, right after `/user/hugs/current`, and drops the old
# `/game/target` / `/game/hugs` rows from the "game" section entirely.

$wb = $excel.ActiveWorkbook
$get = $wb.Worksheets.Item("GET")

# Insert the new `/user/target` row right after `/user/hugs/current` (row 8),
# pushing the remaining /user/* and /friends/* rows down by one.
$get.Rows.Item(9).Insert()
$get.Range("A9").Value = "/user/target"
$get.Range("B9").Value = "id=1234"
$get.Range("C9").Value = "{result:1234}"

# Remove the old `/game/target` (row 19) and `/game/hugs` (row 20) rows
# entirely (those row numbers account for the row just inserted above).
$get.Range("19:20").Delete()

# Leave the selection where the edit happened, and make GET the active tab.
$get.Range("C9").Select() | Out-Null
